$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Select()

$ws.Range("A7").Value = "cell_count"

$ws.Range("A13").Value = "cell_count"
$ws.Range("B13").Value = 29
$ws.Range("C13").Value = 25
$ws.Range("D13").Value = 22
$ws.Range("E13").Value = 26
$ws.Range("F13").Value = 20
$ws.Range("G13").Value = 29
$ws.Range("H13").Value = 29

$ws.Range("I22").Select()
